{"js": "// Replace the placeholder \"<Condition>\" with \"<Assessment of Significance>\"\n// and drop the stray \"_GoBack\" bookmark that Word leaves behind after an\n// edit (it marked the position of the last editing change).\nconst body = context.document.body;\n\nconst results = body.search(\"<Condition>\", { matchCase: true, matchWildcards: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"<Assessment of Significance>\", Word.InsertLocation.replace);\n}\n\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Replace the placeholder \"<Condition>\" with \"<Assessment of Significance>\".\n$d.Content.Find.Execute(\"<Condition>\", $false, $false, $false, $false, $false, $true, 1, $false, \"<Assessment of Significance>\", 2)\n\n# Drop the stray \"_GoBack\" bookmark left over from the last editing session.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
